# Tarea 2 version inicial
# Updates the "DataPruebas" sheet test-data table: refreshes the signup /
# login / cart / checkout sample rows with Cristian Moraga's data set.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataPruebas")

# --- Row 2: CP001_creacion_cta (account creation) ---------------------
$ws.Range("A2").Value = "CP001_creacion_cta"
$ws.Range("B2").Value = "Cristian"
$ws.Range("C2").Value = "Moraga"
$ws.Range("D2").Value = "correo4@gmail.com"
$ws.Range("E2").Value = 977595472
$ws.Range("F2").Value = "clave123"
$ws.Range("G2").Value = "clave123"
$ws.Range("H2").Value = "Congratulations! Your new account has been successfully created!"
$ws.Range("I2").ClearContents()

# --- Row 3: CP002_validar_creacion_cuenta_antigua (duplicate signup) --
$ws.Range("A3").Value = "CP002_validar_creacion_cuenta_antigua"
$ws.Range("B3").Value = "Cristian"
$ws.Range("C3").Value = "Moraga"
$ws.Range("D3").Value = "correo4@gmail.com"
$ws.Range("E3").Value = 977595472
$ws.Range("F3").Value = "clave123"
$ws.Range("G3").Value = "clave123"
$ws.Range("H3").Value = "Warning: E-Mail Address is already registered!"

# --- Row 4: CP003_validar_login_correcto (login) -----------------------
$ws.Range("A4").Value = "CP003_validar_login_correcto"
$ws.Range("C4").Value = 123456
$ws.Range("D4").Value = "My Account"
$ws.Range("E4").Value = "X"
$ws.Range("F4").Value = "X"

# --- Row 5: CP004_agregar_producto_carro (add to cart) -----------------
$ws.Range("A5").Value = "CP004_agregar_producto_carro"
$ws.Range("C5").Value = 123456
$ws.Range("D5").Value = "Success: You have added MacBook to your shopping cart!"

# --- Row 6: CP005_producto_carro_confirm_order (checkout) --------------
$ws.Range("A6").Value = "CP005_producto_carro_confirm_order"
$ws.Range("C6").Value = 123456
$ws.Range("D6").Value = "Your order has been placed!"

# --- Hyperlinks ----------------------------------------------------------
# B2 / B3 already carried mailto hyperlinks before this edit; typing the new
# name over them (above) leaves those links attached, matching how Excel
# keeps a hyperlink when you overwrite the cell text in place.
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:correo4@gmail.com", "", "", "correo4@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:correo4@gmail.com", "", "", "correo4@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:cmoraga.ochoa@gmail.com", "", "", "cmoraga.ochoa@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:cmoraga.ochoa@gmail.com", "", "", "cmoraga.ochoa@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B6"), "mailto:cmoraga.ochoa@gmail.com", "", "", "cmoraga.ochoa@gmail.com")

# Re-apply the workbook's built-in hyperlink style only to the freshly
# linked cells. B2/B3 carried the blue/underlined "Hipervínculo" look from
# the old data; once their text is overwritten in place that formatting is
# no longer applied (only the hyperlink field itself survives), so leave
# them at the default style.
$ws.Range("D2").Style = "Hipervínculo"
$ws.Range("D3").Style = "Hipervínculo"
$ws.Range("B4").Style = "Hipervínculo"
$ws.Range("B5").Style = "Hipervínculo"
$ws.Range("B6").Style = "Hipervínculo"

# --- Column widths (auto-fit sizing after the longer sample data) -------
$ws.Columns.Item(1).ColumnWidth = 21.86
$ws.Columns.Item(2).ColumnWidth = 21.29
$ws.Columns.Item(4).ColumnWidth = 51.86
$ws.Columns.Item(8).ColumnWidth = 60.29

# --- Selection left where the author finished editing --------------------
$ws.Range("D15").Select()
